$d = $word.ActiveDocument

$pairs = @(
    @("420×7=", "333×6="),
    @("965×5=", "576×5="),
    @("330×4=", "568×4="),
    @("459×6=", "647×8="),
    @("258×3=", "173×3="),
    @("335×5=", "362×6="),
    @("470×7=", "537×4="),
    @("462×7=", "524×5="),
    @("354×9=", "704×4="),
    @("336×3=", "967×3="),
    @("656×9=", "300×2="),
    @("511×4=", "302×5="),
    @("495×5=", "992×4="),
    @("139×2=", "234×2="),
    @("411×9=", "877×9="),
    @("744×3=", "960×6="),
    @("888×8=", "403×7="),
    @("263×8=", "252×8="),
    @("447×8=", "172×9="),
    @("390×7=", "703×2="),
    @("177×7=", "321×9="),
    @("860×4=", "156×6="),
    @("909×7=", "967×2="),
    @("656×6=", "783×5="),
    @("332×9=", "104×2=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
